$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.013.12"
$ws.Range("E2").Value = "  -0.75%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.420.45"
$ws.Range("E3").Value = "  -0.50%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "411.04"
$ws.Range("E5").Value = "  +0.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.46"
$ws.Range("E6").Value = "  -3.44%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.642"
$ws.Range("E7").Value = "  +8.27%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.740"
$ws.Range("E9").Value = "  +7.96%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.142"
$ws.Range("E10").Value = "  +15.58%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.78"
$ws.Range("E11").Value = "  +0.93%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000219"
$ws.Range("E12").Value = "  +66.53%  "
$ws.Range("E13").Value = "  +8.35%  "
$ws.Range("E14").Value = "  -0.39%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.954.57"
$ws.Range("E15").Value = "  -0.64%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.26"
$ws.Range("E16").Value = "  +6.69%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.434.36"
$ws.Range("E17").Value = "  -0.41%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.08"
$ws.Range("E18").Value = "  +5.35%  "
$ws.Range("E19").Value = "  +5.93%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "61.973.89"
$ws.Range("E20").Value = "  -0.69%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "445.86"
$ws.Range("E21").Value = "  +41.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "91.78"
$ws.Range("E22").Value = "  +9.04%  "
$ws.Range("E23").Value = "  -0.47%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.12"
$ws.Range("E24").Value = "  +1.29%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.27"
$ws.Range("E25").Value = "  +3.29%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "33.53"
$ws.Range("E26").Value = "  +12.60%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.85"
$ws.Range("E27").Value = "  +7.29%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "4.75"
$ws.Range("E28").Value = "  +0.70%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.62"
$ws.Range("E29").Value = "  +0.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.75"
$ws.Range("E30").Value = "  +0.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "12.02"
$ws.Range("E31").Value = "  +5.44%  "
$ws.Range("E32").Value = "  -0.06%  "
$ws.Range("E33").Value = "  -2.35%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "42.85"
$ws.Range("E34").Value = "  +1.42%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  -0.10%  "
$ws.Range("E36").Value = "  +3.52%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "53.71"
$ws.Range("E37").Value = "  +4.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.998"
$ws.Range("E38").Value = "  +0.09%  "
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.39"
$ws.Range("E39").Value = "  -0.79%  "
$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.93"
$ws.Range("E40").Value = "  -1.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.316"
$ws.Range("E42").Value = "  -0.23%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "141.46"
$ws.Range("E43").Value = "  +2.77%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.24"
$ws.Range("E44").Value = "  +5.26%  "
$ws.Range("E45").Value = "  -0.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.40"
$ws.Range("E46").Value = "  +7.61%  "
$ws.Range("E47").Value = "  -1.22%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.33"
$ws.Range("E48").Value = "  +4.36%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.765.22"
$ws.Range("E49").Value = "  -0.31%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.113.37"
$ws.Range("E50").Value = "  -0.73%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "105.20"
$ws.Range("E51").Value = "  +25.46%  "
